# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.028.99"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "1.744.05"
$ws.Range("E3").Value = "  +1.41%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.38%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.11"
$ws.Range("E5").Value = "  -1.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9997"
$ws.Range("E6").Value = "  -0.47%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4977"
$ws.Range("E7").Value = "  +7.43%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3585"
$ws.Range("E8").Value = "  +4.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.64"
$ws.Range("E9").Value = "  +0.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07276"
$ws.Range("E10").Value = "  -0.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.063"
$ws.Range("E11").Value = "  +1.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9992"
$ws.Range("E12").Value = "  -0.47%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.32"
$ws.Range("E13").Value = "  +2.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.978"
$ws.Range("E14").Value = "  +1.62%  "
$ws.Range("D15").Value = "1.746.72"
$ws.Range("E15").Value = "  +1.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.878"
$ws.Range("E16").Value = "  -0.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "87.07"
$ws.Range("E17").Value = "  -2.78%  "
$ws.Range("E18").Value = "  -0.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06373"
$ws.Range("E19").Value = "  +0.95%  "
$ws.Range("E20").Value = "  -0.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.63"
$ws.Range("E21").Value = "  +0.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.724"
$ws.Range("E22").Value = "  +1.51%  "
$ws.Range("D23").Value = "27.109.91"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.32"
$ws.Range("E24").Value = "  +4.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.039"
$ws.Range("E25").Value = "  -5.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.05"
$ws.Range("E26").Value = "  -1.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.96"
$ws.Range("E27").Value = "  +2.45%  "
$ws.Range("D28").Value = "1.939.12"
$ws.Range("E28").Value = "  +1.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.170"
$ws.Range("E29").Value = "  +0.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.59"
$ws.Range("E30").Value = "  +1.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.052"
$ws.Range("E31").Value = "  +2.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09564"
$ws.Range("E32").Value = "  +5.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.573"
$ws.Range("E33").Value = "  -0.80%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.407"
$ws.Range("E34").Value = "  +1.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02202"
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05897"
$ws.Range("E36").Value = "  +0.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "11.08"
$ws.Range("E37").Value = "  -0.60%  "
$ws.Range("E38").Value = "  +2.31%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.791"
$ws.Range("E39").Value = "  +0.16%  "
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2005"
$ws.Range("E40").Value = "  +0.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6046"
$ws.Range("E41").Value = "  +1.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.114"
$ws.Range("E42").Value = "  -1.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.573"
$ws.Range("E43").Value = "  +0.77%  "
$ws.Range("E44").Value = "  +2.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.595"
$ws.Range("E45").Value = "  -0.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5680"
$ws.Range("E46").Value = "  +0.75%  "
$ws.Range("E47").Value = "  +0.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.868"
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.104"
$ws.Range("E49").Value = "  +1.48%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06685"
$ws.Range("E50").Value = "  +0.29%  "
$ws.Range("B51").Value = "PaxDollar"
$ws.Range("C51").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9995"
$ws.Range("E51").Value = "  -0.56%  "
